$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G1 gets its own (non-shared) formula, matching E1/C1 pattern
$ws.Range("G1").Formula = "=E1^2"

# G2:G16 share one relative formula, mirroring how C2:C16 / E2:E16 are shared
$ws.Range("G2:G16").FormulaR1C1 = "=RC[-2]^2"

# Column G: best-fit width (~11.625 characters, matching Excel's AutoFit result)
$ws.Columns.Item(7).AutoFit() | Out-Null
$ws.Columns.Item(7).ColumnWidth = 10.86

# Update the selected cell to G3
$ws.Range("G3").Select()
